$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 7 ("cap to xor 0x70E040 I xora"), pushing the
# existing rows 7-10 (ypos, ground lvl, speed/pz, follow target) down by one.
$ws.Rows.Item(7).Insert()
$ws.Range("A7").Value = "cap to xor 0x70E040 I xora"

# Insert a new row at row 12 ("server log - najnowsza linia to 0x5C3DC0"),
# right after the "follow target" line (now at row 11), and remove the
# following blank row so that the gap before "mousel last X click..." stays
# a single empty row.
$ws.Rows.Item(12).Insert()
$ws.Range("A12").Value = "server log - najnowsza linia to 0x5C3DC0"
$ws.Rows.Item(13).Delete()

# Fix the selection to match the recorded state after the edit.
$ws.Range("H23").Select()
